$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'56.352.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.15%  '
$ws.Range("D3").Value = "'2.970.84"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -5.27%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'495.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.44%  '
$ws.Range("D6").Value = "'134.53"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.45%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = "'2.967.93"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -5.31%  '
$ws.Range("E9").Value = '  -3.91%  '
$ws.Range("E10").Value = '  -0.86%  '
$ws.Range("E11").Value = '  -3.71%  '
$ws.Range("E12").Value = '  -7.37%  '
$ws.Range("E13").Value = '  -0.77%  '
$ws.Range("D14").Value = "'3.480.72"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.20%  '
$ws.Range("D15").Value = "'25.04"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.60%  '
$ws.Range("D16").Value = "'56.380.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.06%  '
$ws.Range("D17").Value = "'2.972.16"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.03%  '
$ws.Range("E18").Value = '  -4.56%  '
$ws.Range("E19").Value = '  +1.02%  '
$ws.Range("D20").Value = "'12.35"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.26%  '
$ws.Range("E21").Value = '  -1.92%  '
$ws.Range("D22").Value = "'325.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -5.31%  '
$ws.Range("E23").Value = '  -0.04%  '
$ws.Range("E24").Value = '  -7.82%  '
$ws.Range("E25").Value = '  -9.49%  '
$ws.Range("E26").Value = '  -1.08%  '
$ws.Range("E27").Value = '  -5.25%  '
$ws.Range("D28").Value = "'0.0₃0896"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.94%  '
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("D30").Value = "'6.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.27%  '
$ws.Range("D31").Value = "'6.71"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.49%  '
$ws.Range("E32").Value = '  -5.08%  '
$ws.Range("E33").Value = '  -6.68%  '
$ws.Range("D34").Value = "'20.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.73%  '
$ws.Range("D35").Value = "'152.79"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.93%  '
$ws.Range("E36").Value = '  -8.04%  '
$ws.Range("E37").Value = '  -6.69%  '
$ws.Range("D38").Value = "'5.59"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -10.30%  '
$ws.Range("D39").Value = "'0.0669"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.17%  '
$ws.Range("D40").Value = "'23.13"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.85%  '
$ws.Range("D41").Value = "'3.005.23"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'36.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -9.63%  '
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("E44").Value = '  -7.79%  '
$ws.Range("D45").Value = "'0.992"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -9.04%  '
$ws.Range("D46").Value = "'1.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.13%  '
$ws.Range("D47").Value = "'2.206.52"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.27%  '
$ws.Range("E48").Value = '  -8.92%  '
$ws.Range("D49").Value = "'1.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.01%  '
$ws.Range("E50").Value = '  +1.48%  '
$ws.Range("D51").Value = "'5.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.95%  '
